$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 6071
$ws.Range("F6").Value = 15
$ws.Range("F10").Value = 705
$ws.Range("F11").Value = 1958
$ws.Range("F12").Value = 1958
$ws.Range("F14").Value = 1692
$ws.Range("F15").Value = 586
$ws.Range("F16").Value = 226
$ws.Range("F18").Value = 4806
$ws.Range("G18").Value = 80
$ws.Range("F19").Value = 131
$ws.Range("F22").Value = 3373
$ws.Range("F23").Value = 840
$ws.Range("F25").Value = 61
$ws.Range("F27").Value = 2377
$ws.Range("F33").Value = 1259
$ws.Range("F35").Value = 38
$ws.Range("F38").Value = 1339
$ws.Range("F39").Value = 1315

# --- Sheet: 演出 (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F14").Value = 105
$ws.Range("F19").Value = 134
$ws.Range("F20").Value = 319
$ws.Range("F21").Value = 248

# --- Sheet: 全部类型 (All types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 6071
$ws.Range("F14").Value = 15
$ws.Range("F21").Value = 1959
$ws.Range("F23").Value = 1692
$ws.Range("F24").Value = 105
$ws.Range("F25").Value = 586
$ws.Range("F26").Value = 226
$ws.Range("F28").Value = 4806
$ws.Range("F31").Value = 3373
$ws.Range("F33").Value = 61
$ws.Range("F36").Value = 2377
$ws.Range("F40").Value = 1259
$ws.Range("F41").Value = 134
$ws.Range("F42").Value = 248
$ws.Range("F45").Value = 38
$ws.Range("F48").Value = 1339
